# Update workbook data to the "output generated at 456a3b4" snapshot.
# The workbook has 4 sheets: 展览 (Exhibition), 演出 (Performance),
# 本地生活 (Local life), 全部类型 (All types).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet: 展览 (Exhibition)
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("展览")

$ws.Range("F4").Value = 595
$ws.Range("F5").Value = 9259
$ws.Range("F7").Value = 11981
$ws.Range("F8").Value = 11981

# Row 9 now holds what used to be row 9's "next" entry (内田秀), row 10
# shifts to what was row 9 (小林爱香), row 11 shifts to what was row 10
# (青山渚), replacing the old row 11 (广播剧《西东》).
$ws.Range("C9").Value = "北京·人气声优 内田秀 专场活动"
$ws.Range("E9").Value = "2024.10.02 13:55-10.02 17:10"
$ws.Range("F9").Value = 130
$ws.Range("H9").Value = "https://show.bilibili.com/platform/detail.html?id=91678"
$ws.Range("I9").Value = "//i0.hdslb.com/bfs/openplatform/202409/0aUkHD511725260741169.png"

$ws.Range("C10").Value = "北京·人气声优 小林爱香 专场活动"
$ws.Range("E10").Value = "2024.10.02 12:50-10.02 16:40"
$ws.Range("F10").Value = 228
$ws.Range("H10").Value = "https://show.bilibili.com/platform/detail.html?id=91117"
$ws.Range("I10").Value = "//i2.hdslb.com/bfs/openplatform/202408/nuqS5Gd11724309352207.png"

$ws.Range("C11").Value = "北京·人气声优 青山渚 专场活动"
$ws.Range("E11").Value = "2024.10.02 11:50-10.02 15:40"
$ws.Range("F11").Value = 309
$ws.Range("H11").Value = "https://show.bilibili.com/platform/detail.html?id=91249"
$ws.Range("I11").Value = "//i2.hdslb.com/bfs/openplatform/202408/xHqpdFa41724641733192.png"

$ws.Range("F13").Value = 127
$ws.Range("F15").Value = 447
$ws.Range("F17").Value = 2076
$ws.Range("F18").Value = 831
$ws.Range("F19").Value = 790
$ws.Range("F20").Value = 395
$ws.Range("F21").Value = 50
$ws.Range("F22").Value = 414
$ws.Range("F23").Value = 313
$ws.Range("F25").Value = 669
$ws.Range("F26").Value = 22
$ws.Range("F27").Value = 1567
$ws.Range("F29").Value = 24
$ws.Range("F30").Value = 21
$ws.Range("F31").Value = 56
$ws.Range("F33").Value = 1445
$ws.Range("F34").Value = 9
$ws.Range("F35").Value = 505
$ws.Range("F36").Value = 342
$ws.Range("F37").Value = 551
$ws.Range("F38").Value = 399
$ws.Range("F39").Value = 2216
$ws.Range("F42").Value = 158
$ws.Range("F43").Value = 580
$ws.Range("F44").Value = 451
$ws.Range("F45").Value = 167
$ws.Range("F46").Value = 889
$ws.Range("F47").Value = 684
$ws.Range("F49").Value = 326
$ws.Range("F50").Value = 295

# ---------------------------------------------------------------
# Sheet: 演出 (Performance)
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("演出")

$ws.Range("F6").Value = 75
$ws.Range("F12").Value = 48
$ws.Range("F22").Value = 79
$ws.Range("F24").Value = 72

# ---------------------------------------------------------------
# Sheet: 本地生活 (Local life)
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("本地生活")

$ws.Range("F4").Value = 369
$ws.Range("F5").Value = 238
$ws.Range("F6").Value = 277

# ---------------------------------------------------------------
# Sheet: 全部类型 (All types)
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("全部类型")

$ws.Range("F6").Value = 369
$ws.Range("F7").Value = 238
$ws.Range("F8").Value = 595
$ws.Range("F9").Value = 9259
$ws.Range("F16").Value = 447
$ws.Range("F17").Value = 50
$ws.Range("F18").Value = 414
$ws.Range("F19").Value = 313
$ws.Range("F21").Value = 669
$ws.Range("F22").Value = 22
$ws.Range("F23").Value = 277
$ws.Range("F24").Value = 1567
$ws.Range("F26").Value = 48
$ws.Range("F28").Value = 56
$ws.Range("F32").Value = 1445
$ws.Range("F34").Value = 9
$ws.Range("F35").Value = 505
$ws.Range("F36").Value = 551
$ws.Range("F37").Value = 399
$ws.Range("F39").Value = 2216
$ws.Range("F41").Value = 158
$ws.Range("F42").Value = 580
$ws.Range("F43").Value = 451
$ws.Range("F44").Value = 167
$ws.Range("F45").Value = 889
$ws.Range("F47").Value = 79
$ws.Range("F49").Value = 684
